# SG2042EVB(5): fixbug: pioneer system poweroff failed
# On the "power" worksheet, insert a new row above the existing
# EN_VDD_3V3 row (row 2) for the MILKV_ATX_CTL control signal, shifting
# all subsequent rows (including the new SYS_RST_DEASSERT row) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("power")

# Insert a new blank row at row 2; existing rows 2-22 shift down to 3-23.
$ws.Rows.Item(2).Insert()

# Make sure the Delay column keeps its original "number stored as text"
# formatting before writing the value into it.
$ws.Range("A2:D2").NumberFormat = "@"

$ws.Range("A2").Value = "MILKV_ATX_CTL"
$ws.Range("B2").Value = "MILKV_ATX_CTL"
$ws.Range("C2").Value = "FUNCTION"
$ws.Range("D2").Value = "1000"
